$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1022.32434
$ws.Range("I15").Value = 1022.32434
$ws.Range("K15").Value = 3066.97302
$ws.Range("M15").Value = -2897.97302
$ws.Range("H40").Value = 5613.5557
$ws.Range("I40").Value = 2097.5
$ws.Range("K40").Value = 2097.5
$ws.Range("M40").Value = -1922.5
$ws.Range("H112").Value = 2909.0857
$ws.Range("J112").Value = 2965.853
$ws.Range("L112").Value = 8897.559000000001
$ws.Range("N112").Value = -11113.559
$ws.Range("H125").Value = 6731.875
$ws.Range("I125").Value = 6795.7856
$ws.Range("J125").Value = 6642.4
$ws.Range("K125").Value = 61162.0704
$ws.Range("L125").Value = 59781.6
$ws.Range("M125").Value = -58702.0704
$ws.Range("N125").Value = -64701.6
$ws.Range("H131").Value = 2571.9285
$ws.Range("I131").Value = 1869.8182
$ws.Range("K131").Value = 5609.4546
$ws.Range("M131").Value = -569.4546
$ws.Range("H132").Value = 1825.6207
$ws.Range("I132").Value = 1037.72
$ws.Range("K132").Value = 3113.16
$ws.Range("M132").Value = -583.1599999999999
$ws.Range("H133").Value = 61365.453
$ws.Range("J133").Value = 61365.453
$ws.Range("L133").Value = 61365.453
$ws.Range("N133").Value = -71485.45300000001
$ws.Range("H137").Value = 3971.1353
$ws.Range("J137").Value = 7759.357
$ws.Range("L137").Value = 23278.071
$ws.Range("N137").Value = -28378.071

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 20654.572
$ws.Range("I2").Value = 23812.021
$ws.Range("K2").Value = 23812.021
$ws.Range("M2").Value = -23699.021
$ws.Range("H45").Value = 6460.3335
$ws.Range("I45").Value = 3890.5
$ws.Range("J45").Value = 11600
$ws.Range("K45").Value = 3890.5
$ws.Range("L45").Value = 11600
$ws.Range("M45").Value = -3513.5
$ws.Range("N45").Value = -12354
$ws.Range("H116").Value = 20654.572
$ws.Range("I116").Value = 23812.021
$ws.Range("K116").Value = 23812.021
$ws.Range("M116").Value = -21518.021
$ws.Range("H118").Value = 100000
$ws.Range("J118").Value = 100000
$ws.Range("L118").Value = 100000
$ws.Range("N118").Value = -103314
$ws.Range("H132").Value = 6769.1714
$ws.Range("I132").Value = 3658.8333
$ws.Range("K132").Value = 10976.4999
$ws.Range("M132").Value = -8446.499899999999

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 20654.572
$ws.Range("I3").Value = 23812.021
$ws.Range("K3").Value = 23812.021
$ws.Range("M3").Value = -23698.021
$ws.Range("H86").Value = 1419.8182
$ws.Range("I86").Value = 1095.3334
$ws.Range("K86").Value = 1095.3334
$ws.Range("M86").Value = 27.66660000000002
$ws.Range("H89").Value = 1419.8182
$ws.Range("I89").Value = 1095.3334
$ws.Range("K89").Value = 5476.666999999999
$ws.Range("M89").Value = 139.3330000000005
$ws.Range("H94").Value = 1122.2222
$ws.Range("I94").Value = 1122.2222
$ws.Range("K94").Value = 1122.2222
$ws.Range("M94").Value = -671.2221999999999
$ws.Range("H107").Value = 1534.8334
$ws.Range("I107").Value = 1302.25
$ws.Range("K107").Value = 1302.25
$ws.Range("M107").Value = 617.75

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4023
$ws.Range("I31").Value = 1862.5
$ws.Range("J31").Value = 7479.8
$ws.Range("K31").Value = 1862.5
$ws.Range("L31").Value = 7479.8
$ws.Range("M31").Value = -1567.5
$ws.Range("N31").Value = -8069.8
$ws.Range("H34").Value = 4023
$ws.Range("I34").Value = 1862.5
$ws.Range("J34").Value = 7479.8
$ws.Range("K34").Value = 1862.5
$ws.Range("L34").Value = 7479.8
$ws.Range("M34").Value = -1660.5
$ws.Range("N34").Value = -7883.8
$ws.Range("H104").Value = 46500
$ws.Range("J104").Value = 46500
$ws.Range("L104").Value = 46500
$ws.Range("N104").Value = -51742
$ws.Range("H107").Value = 4433
$ws.Range("I107").Value = 2999
$ws.Range("K107").Value = 2999
$ws.Range("M107").Value = -1079
$ws.Range("H132").Value = 3361.5483
$ws.Range("I132").Value = 2848.3845
$ws.Range("K132").Value = 8545.1535
$ws.Range("M132").Value = -6015.1535
$ws.Range("H134").Value = 3153.75
$ws.Range("J134").Value = 3964.5293
$ws.Range("L134").Value = 11893.5879
$ws.Range("N134").Value = -16963.5879

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 3640484.2
$ws.Range("I4").Value = 5779069.5
$ws.Range("K4").Value = 17337208.5
$ws.Range("M4").Value = -17337096.5
$ws.Range("H5").Value = 1939
$ws.Range("I5").Value = 1932
$ws.Range("J5").Value = 1949.5
$ws.Range("K5").Value = 5796
$ws.Range("L5").Value = 5848.5
$ws.Range("M5").Value = -5684
$ws.Range("N5").Value = -6072.5
$ws.Range("H57").Value = 0
$ws.Range("I57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("M57").ClearContents()
$ws.Range("H101").Value = 13811.6
$ws.Range("J101").Value = 13811.6
$ws.Range("L101").Value = 41434.8
$ws.Range("N101").Value = -46302.8
$ws.Range("H133").Value = 2500
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 2500
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 7500
$ws.Range("M133").ClearContents()
$ws.Range("N133").Value = -17620
$ws.Range("H134").Value = 1799
$ws.Range("I134").Value = 1799
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 5397
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -327
$ws.Range("N134").ClearContents()
$ws.Range("H135").Value = 1939
$ws.Range("I135").Value = 1932
$ws.Range("J135").Value = 1949.5
$ws.Range("K135").Value = 17388
$ws.Range("L135").Value = 17545.5
$ws.Range("M135").Value = -14853
$ws.Range("N135").Value = -22615.5
$ws.Range("H138").Value = 2255.1428
$ws.Range("I138").Value = 957.4
$ws.Range("J138").Value = 5499.5
$ws.Range("K138").Value = 2872.2
$ws.Range("L138").Value = 16498.5
$ws.Range("M138").Value = 2267.8
$ws.Range("N138").Value = -26778.5
$ws.Range("H139").Value = 6371.7
$ws.Range("I139").Value = 4343.4
$ws.Range("J139").Value = 8400
$ws.Range("K139").Value = 13030.2
$ws.Range("L139").Value = 25200
$ws.Range("M139").Value = -7890.199999999999
$ws.Range("N139").Value = -35480
$ws.Range("H140").Value = 3057.25
$ws.Range("I140").Value = 3057.25
$ws.Range("K140").Value = 9171.75
$ws.Range("M140").Value = -3991.75
$ws.Range("H141").Value = 8140.3335
$ws.Range("J141").Value = 9500
$ws.Range("L141").Value = 28500
$ws.Range("N141").Value = -38860

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 128.11765
$ws.Range("I2").Value = 158.53847
$ws.Range("K2").Value = 158.53847
$ws.Range("M2").Value = -45.53846999999999
$ws.Range("H70").Value = 23814872
$ws.Range("I70").Value = 4575.7144
$ws.Range("J70").Value = 47625170
$ws.Range("K70").Value = 4575.7144
$ws.Range("L70").Value = 47625170
$ws.Range("M70").Value = -4305.7144
$ws.Range("N70").Value = -47625710
$ws.Range("H73").Value = 23814872
$ws.Range("I73").Value = 4575.7144
$ws.Range("J73").Value = 47625170
$ws.Range("K73").Value = 4575.7144
$ws.Range("L73").Value = 47625170
$ws.Range("M73").Value = -3639.7144
$ws.Range("N73").Value = -47627042
$ws.Range("H102").Value = 3824.375
$ws.Range("I102").Value = 3729.9312
$ws.Range("K102").Value = 3729.9312
$ws.Range("M102").Value = -2107.9312
$ws.Range("H132").Value = 420248.97
$ws.Range("I132").Value = 479809.06
$ws.Range("J132").Value = 3328.3333
$ws.Range("K132").Value = 1439427.18
$ws.Range("L132").Value = 9984.999899999999
$ws.Range("M132").Value = -1436897.18
$ws.Range("N132").Value = -15044.9999

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 631568.25
$ws.Range("J7").Value = 5021.778
$ws.Range("L7").Value = 5021.778
$ws.Range("N7").Value = -5245.778
$ws.Range("H40").Value = 672832.9
$ws.Range("I40").Value = 1004849.5
$ws.Range("J40").Value = 8799.6
$ws.Range("K40").Value = 1004849.5
$ws.Range("L40").Value = 8799.6
$ws.Range("M40").Value = -1004713.5
$ws.Range("N40").Value = -9071.6
$ws.Range("H46").Value = 3410.1614
$ws.Range("I46").Value = 2813.4348
$ws.Range("K46").Value = 2813.4348
$ws.Range("M46").Value = -2625.4348
$ws.Range("H122").Value = 3503501.8
$ws.Range("I122").Value = 2505250.2
$ws.Range("K122").Value = 7515750.600000001
$ws.Range("M122").Value = -7513300.600000001
$ws.Range("H126").Value = 631568.25
$ws.Range("J126").Value = 5021.778
$ws.Range("L126").Value = 15065.334
$ws.Range("N126").Value = -20005.334
$ws.Range("H132").Value = 4776.5557
$ws.Range("I132").Value = 3461.6365
$ws.Range("K132").Value = 10384.9095
$ws.Range("M132").Value = -7854.9095
$ws.Range("H141").Value = 50000
$ws.Range("J141").Value = 50000
$ws.Range("L141").Value = 50000
$ws.Range("N141").Value = -60360

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 144585.28
$ws.Range("I96").Value = 144585.28
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 144585.28
$ws.Range("L96").Value = 0
$ws.Range("M96").Value = -143212.28
$ws.Range("N96").ClearContents()
$ws.Range("H132").Value = 3329.5715
$ws.Range("I132").Value = 2583.5881
$ws.Range("J132").Value = 6500
$ws.Range("K132").Value = 7750.7643
$ws.Range("L132").Value = 19500
$ws.Range("M132").Value = -5220.7643
$ws.Range("N132").Value = -24560
$ws.Range("H136").Value = 1671417.1
$ws.Range("I136").Value = 2004100.6
$ws.Range("K136").Value = 6012301.800000001
$ws.Range("M136").Value = -6009751.800000001
